$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.119.27"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "2.050.21"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "251.47"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "0.666"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "56.41"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("D9").Value = "61.55"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("D10").Value = "0.385"
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("D11").Value = "0.0784"
$ws.Range("E11").Value = "  +3.73%  "
$ws.Range("E12").Value = "  +1.64%  "
$ws.Range("D13").Value = "16.49"
$ws.Range("E13").Value = "  +6.10%  "
$ws.Range("D14").Value = "2.341.97"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").Value = "0.803"
$ws.Range("E15").Value = "  -4.84%  "
$ws.Range("D16").Value = "5.56"
$ws.Range("E16").Value = "  +4.82%  "
$ws.Range("D17").Value = "2.044.59"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").Value = "37.060.75"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "16.74"
$ws.Range("E19").Value = "  +13.36%  "
$ws.Range("D20").Value = "75.14"
$ws.Range("E20").Value = "  +2.61%  "
$ws.Range("D21").Value = "0.0₃0908"
$ws.Range("E21").Value = "  +6.63%  "
$ws.Range("D22").Value = "5.39"
$ws.Range("E22").Value = "  +2.09%  "
$ws.Range("D23").Value = "237.17"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("E25").Value = "  -4.00%  "
$ws.Range("E26").Value = "  +17.45%  "
$ws.Range("D27").Value = "169.07"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").Value = "9.25"
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").Value = "20.17"
$ws.Range("E29").Value = "  -4.33%  "
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("E31").Value = "  +5.01%  "
$ws.Range("D32").Value = "4.72"
$ws.Range("E32").Value = "  +3.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0620"
$ws.Range("E33").Value = "  -1.44%  "
$ws.Range("D34").Value = "4.45"
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("D35").Value = "0.0885"
$ws.Range("E35").Value = "  -3.07%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "2.24"
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("E38").Value = "  -4.48%  "
$ws.Range("D39").Value = "0.108"
$ws.Range("E39").Value = "  +8.08%  "
$ws.Range("D40").Value = "1.35"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("D41").Value = "17.74"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").Value = "0.0224"
$ws.Range("E42").Value = "  -2.30%  "
$ws.Range("D43").Value = "1.14"
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("D44").Value = "96.78"
$ws.Range("E44").Value = "  -3.33%  "
$ws.Range("B45").Value = "THORChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D45").Value = "4.76"
$ws.Range("E45").Value = "  +18.25%  "
$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").Value = "2.84"
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("D47").Value = "2.49"
$ws.Range("E47").Value = "  +3.32%  "
$ws.Range("D48").Value = "1.285.12"
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.90"
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("D50").Value = "6.76"
$ws.Range("E50").Value = "  -3.84%  "
$ws.Range("D51").Value = "2.229.63"
$ws.Range("E51").Value = "  -1.61%  "
